# Updated symbol list on Mon Jan 23 14:42:33 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2"  = "304.89"; "E2"  = "-0.04%"
    "D3"  = "35.52";  "E3"  = "-4.19%"
    "D4"  = "5.066";  "E4"  = "1.10%"
    "D5"  = "0.07881"; "E5" = "-0.21%"
    "D6"  = "2.112";  "E6"  = "-4.36%"
    "D7"  = "4.127";  "E7"  = "2.77%"
    "D8"  = "7.903";  "E8"  = "-1.53%"
    "D9"  = "0.9219"; "E9"  = "0.08%"
    "D10" = "0.09646"; "E10" = "-0.26%"
    "E11" = "-2.49%"
    "D12" = "0.08643"; "E12" = "0.17%"
    "D13" = "0.03558"; "E13" = "-3.43%"
    "D14" = "0.09929"; "E14" = "-0.61%"
    "D16" = "0.005653"; "E16" = "0.40%"
    "E17" = "-0.14%"
    "D18" = "2.640"; "E18" = "17.39%"
    "D19" = "0.3371"; "E19" = "-1.26%"
    "D20" = "0.1341"
    "D21" = "5.159"; "E21" = "8.52%"
    "D22" = "0.2214"; "E22" = "0.70%"
    "D23" = "0.04523"; "E23" = "-0.73%"
    "E24" = "-0.08%"
    "D25" = "0.004852"; "E25" = "8.48%"
    "E26" = "-7.05%"
    "D27" = "0.0004765"; "E27" = "0.29%"
    "D39" = "0.01833"; "E39" = "-0.53%"
    "D40" = "0.04726"; "E40" = "-0.80%"
    "D41" = "0.007895"; "E41" = "-2.95%"
    "D42" = "0.1390"; "E42" = "-0.68%"
    "D43" = "0.007766"; "E43" = "2.79%"
    "D44" = "0.002223"; "E44" = "-0.35%"
    "D45" = "0.01115"; "E45" = "10.90%"
    "D46" = "0.00006374"; "E46" = "1.70%"
    "D47" = "0.00000000753"; "E47" = "0.31%"
    "E48" = "0.24%"
    "D49" = "50.61"; "E49" = "35.00%"
    "D50" = "0.001906"; "E50" = "10.82%"
    "D51" = "0.00002107"; "E51" = "0.31%"
}

foreach ($key in $updates.Keys) {
    $cell = $ws.Range($key)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$key]
}
